$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# New "web" (column X) command added: rightClick(locator)
# Alphabetical insertion at row 79 pushes saveAllWindowIds(var)..waitForTitle(text)
# down by one row (old X79:X122 -> new X80:X123).
#
# New "xml" (column AC) commands added: insertAfter/insertBefore/replaceIn
# Alphabetical insertion shifts AC13:AC18 and extends the list to AC19:AC21.
#
# Rather than physically inserting rows (which would also shift unrelated
# columns like "desktop" in column G), every affected cell in columns X and
# AC is written directly with its final, alphabetically-correct value.
# ---------------------------------------------------------------------------

$updates = @(
    @{Addr="X79"; Value="rightClick(locator)"},
    @{Addr="X80"; Value="saveAllWindowIds(var)"},
    @{Addr="X81"; Value="saveAllWindowNames(var)"},
    @{Addr="X82"; Value="saveAttribute(var,locator,attrName)"},
    @{Addr="X83"; Value="saveAttributeList(var,locator,attrName)"},
    @{Addr="X84"; Value="saveCount(var,locator)"},
    @{Addr="X85"; Value="saveDivsAsCsv(headers,rows,cells,nextPage,file)"},
    @{Addr="X86"; Value="saveElement(var,locator)"},
    @{Addr="X87"; Value="saveElements(var,locator)"},
    @{Addr="X88"; Value="saveLocalStorage(var,key)"},
    @{Addr="X89"; Value="saveLocation(var)"},
    @{Addr="X90"; Value="savePageAs(var,sessionIdName,url)"},
    @{Addr="X91"; Value="savePageAsFile(sessionIdName,url,file)"},
    @{Addr="X92"; Value="saveTableAsCsv(locator,nextPageLocator,file)"},
    @{Addr="X93"; Value="saveText(var,locator)"},
    @{Addr="X94"; Value="saveTextArray(var,locator)"},
    @{Addr="X95"; Value="saveTextSubstringAfter(var,locator,delim)"},
    @{Addr="X96"; Value="saveTextSubstringBefore(var,locator,delim)"},
    @{Addr="X97"; Value="saveTextSubstringBetween(var,locator,start,end)"},
    @{Addr="X98"; Value="saveValue(var,locator)"},
    @{Addr="X99"; Value="scrollLeft(locator,pixel)"},
    @{Addr="X100"; Value="scrollRight(locator,pixel)"},
    @{Addr="X101"; Value="scrollTo(locator)"},
    @{Addr="X102"; Value="select(locator,text)"},
    @{Addr="X103"; Value="selectFrame(locator)"},
    @{Addr="X104"; Value="selectMulti(locator,array)"},
    @{Addr="X105"; Value="selectMultiOptions(locator)"},
    @{Addr="X106"; Value="selectText(locator)"},
    @{Addr="X107"; Value="selectWindow(winId)"},
    @{Addr="X108"; Value="selectWindowAndWait(winId,waitMs)"},
    @{Addr="X109"; Value="selectWindowByIndex(index)"},
    @{Addr="X110"; Value="selectWindowByIndexAndWait(index,waitMs)"},
    @{Addr="X111"; Value="toggleSelections(locator)"},
    @{Addr="X112"; Value="type(locator,value)"},
    @{Addr="X113"; Value="typeKeys(locator,value)"},
    @{Addr="X114"; Value="uncheckAll(locator)"},
    @{Addr="X115"; Value="unselectAllText()"},
    @{Addr="X116"; Value="upload(fieldLocator,file)"},
    @{Addr="X117"; Value="verifyContainText(locator,text)"},
    @{Addr="X118"; Value="verifyText(locator,text)"},
    @{Addr="X119"; Value="wait(waitMs)"},
    @{Addr="X120"; Value="waitForElementPresent(locator)"},
    @{Addr="X121"; Value="waitForPopUp(winId,waitMs)"},
    @{Addr="X122"; Value="waitForTextPresent(text)"},
    @{Addr="X123"; Value="waitForTitle(text)"},
    @{Addr="AC13"; Value="insertAfter(xml,xpath,content,var)"},
    @{Addr="AC14"; Value="insertBefore(xml,xpath,content,var)"},
    @{Addr="AC15"; Value="minify(xml,var)"},
    @{Addr="AC16"; Value="prepend(xml,xpath,content,var)"},
    @{Addr="AC17"; Value="replace(xml,xpath,content,var)"},
    @{Addr="AC18"; Value="replaceIn(xml,xpath,content,var)"},
    @{Addr="AC19"; Value="storeCount(xml,xpath,var)"},
    @{Addr="AC20"; Value="storeValue(xml,xpath,var)"},
    @{Addr="AC21"; Value="storeValues(xml,xpath,var)"},
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# Keep the named ranges in sync with the now-longer lists.
$wb.Names.Item("web").RefersTo = "='#system'!`$X`$2:`$X`$123"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AC`$2:`$AC`$21"
